# Saldo_guide.xlsx update: refresh report date (2024-09-06 -> 2024-09-09)
# and adjust Saldo Previsto / Vl. Total for a subset of accounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/report tab to reflect the new extraction timestamp.
$ws.Name = "IClientBalance-20240909-090518-"

# New values (Saldo Previsto / Vl. Total) for the rows whose balances changed.
$changes = @{
    6   = 1068.1500000000001
    8   = 11298.18
    15  = 2191.81
    17  = 1305.1099999999999
    43  = 1089.53
    49  = 962.09
    51  = 12512.07
    52  = 1534.93
    57  = 1949.77
    60  = 11023.68
    97  = 1204.8699999999999
    101 = 93.97
    102 = 1921.1
    105 = 680.16
    107 = 1349.78
    109 = 367.41
    110 = 3289.48
    112 = 0.01
    120 = 987.05
    138 = 15880.22
    143 = 17734.25
    230 = 25939.7
    245 = 3841.77
    255 = 321.7
}

$lastRow = 274

for ($row = 2; $row -le $lastRow; $row++) {
    # Column G = Dt. Referencia: every row moves from 2024-09-06 to 2024-09-09.
    $ws.Cells.Item($row, 7).Value = 45544

    if ($changes.ContainsKey($row)) {
        $newVal = $changes[$row]
        # Column E = Saldo Previsto, Column H = Vl. Total
        $ws.Cells.Item($row, 5).Value = $newVal
        $ws.Cells.Item($row, 8).Value = $newVal
    }
}
